$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.940.37'
$ws.Range("E2").Value = '  -1.20%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.638.40'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.44'
$ws.Range("E5").Value = '  -0.81%  '
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("E8").Value = '  -0.71%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.63'
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.03'
$ws.Range("E12").Value = '  -0.48%  '
$ws.Range("E13").Value = '  -1.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.632.52'
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("E15").Value = '  -0.95%  '
$ws.Range("E16").Value = '  -0.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.00'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.962.89'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.98'
$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("E23").Value = '  -1.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.23'
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("E25").Value = '  +1.30%  '
$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.28%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.130'
$ws.Range("E27").Value = '  +4.18%  '
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.25'
$ws.Range("E30").Value = '  -0.66%  '
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("E33").Value = '  -0.38%  '
$ws.Range("E34").Value = '  -4.45%  '
$ws.Range("E35").Value = '  +1.74%  '
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.138.37'
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("E41").Value = '  +0.34%  '
$ws.Range("E42").Value = '  -3.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.32'
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.798'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.774.93'
$ws.Range("E45").Value = '  -0.42%  '
$ws.Range("E46").Value = '  +2.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.67'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  +2.97%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.68'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -0.75%  '
